# Update the RAD "SSNmoreThan9Error" test-data timestamps in column B
# (rows 2-5) to reflect the latest test execution run, as captured by
# the "Added RAD BeforePayments Test Cases and Test Data" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = "Mon Sep 11 14:04:35 EDT 2023"
$ws.Range("B3").Value = "Mon Sep 11 14:04:48 EDT 2023"
$ws.Range("B4").Value = "Mon Sep 11 14:05:00 EDT 2023"
$ws.Range("B5").Value = "Mon Sep 11 14:05:13 EDT 2023"
